$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---------------------------------------------------------------------
# Row 1 used to just be a duplicate of row 2's data. It becomes a real
# header row (like the other sheets) and gains the same 7 metadata
# columns (G:M) that every other "fact" sheet already carries.
# ---------------------------------------------------------------------

# Extend the bold/bordered header formatting from B1:F1 into G1:M1.
$ws.Range("B1:F1").Copy() | Out-Null
$ws.Range("G1:M1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# ---------------------------------------------------------------------
# Rows 2-12 keep their existing bank/deposit_type/currency/owner/total
# (columns B:F) values untouched; only the new metadata columns G:M are
# populated.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 12; $r++) {
    # Extend that row's existing (unbordered/non-bold) formatting from
    # B:F across the new G:M columns.
    $ws.Range("B$r`:F$r").Copy() | Out-Null
    $ws.Range("G$r`:M$r").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0

    $idx = $ws.Range("A$r").Value2

    $ws.Range("G$r").Value = "deposit"
    $ws.Range("H$r").Value = "normal"

    # "2012-04-24" looks like a date, so Excel's smart-entry would parse
    # it into a date serial if assigned directly. Force text entry, then
    # re-apply the row's plain formatting so no stray number format is
    # left behind on the cell.
    $ws.Range("I$r").NumberFormat = "@"
    $ws.Range("I$r").Value = "2012-04-24"

    $ws.Range("J$r").Value = "李俊俋"
    $ws.Range("K$r").Value = 1738
    $ws.Range("L$r").Value = "tmp16861"
    $ws.Range("M$r").Value = $idx

    $ws.Range("B$r`:F$r").Copy() | Out-Null
    $ws.Range("G$r`:M$r").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0
}
